$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column A's width (in ColumnWidth units) so the new column can match it.
$colAWidth = $ws.Columns("A:A").ColumnWidth

# Insert a new column before column B (shifts VIN/Effective Date/Expiration Date/Annual GWP right).
$ws.Columns("B:B").Insert()

# Give the new "State" column the same width as "Company Name" (col A), matching Excel's
# behavior of carrying the left-neighbor's formatting into a freshly inserted column.
$ws.Columns("B:B").ColumnWidth = $colAWidth

# Header
$ws.Range("B1").Value2 = "State"

# Company -> State mapping
$companyState = @{
    "Bigtime Trucking"   = "IL"
    "Good Truck Ltd."    = "TN"
    "IJ Asset Delivery"  = "IL"
    "ZoomZoom Go"        = "TN"
    "Cannery Delivery"   = "TN"
    "Python Express"     = "IL"
    "R.R. Fast"          = "TN"
}

for ($r = 2; $r -le 51; $r++) {
    $company = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value2 = $companyState[$company]
}

# Match the final selection state left behind in the saved file.
$ws.Range("B52").Select()
